$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "29.189.20"
$ws.Range("E2").Value = "  -0.54%  "

$ws.Range("D3").Value = "1.827.32"
$ws.Range("E3").Value = "  -0.77%  "

Set-TextCell "D4" "0.9987"
$ws.Range("E4").Value = "  -0.01%  "

Set-TextCell "D5" "236.46"
$ws.Range("E5").Value = "  -1.44%  "

$ws.Range("E6").Value = "  -3.52%  "

Set-TextCell "D7" "1.000"
$ws.Range("E7").Value = "  +0.02%  "

Set-TextCell "D8" "0.07087"
$ws.Range("E8").Value = "  -4.77%  "

Set-TextCell "D9" "0.2811"
$ws.Range("E9").Value = "  -2.83%  "

Set-TextCell "D10" "23.76"
$ws.Range("E10").Value = "  -4.85%  "

Set-TextCell "D11" "0.07662"
$ws.Range("E11").Value = "  -0.83%  "

$ws.Range("D12").Value = "1.818.67"
$ws.Range("E12").Value = "  -1.20%  "

Set-TextCell "D13" "4.817"
$ws.Range("E13").Value = "  -3.10%  "

$ws.Range("E14").Value = "  -2.42%  "

Set-TextCell "D15" "0.6347"
$ws.Range("E15").Value = "  -6.08%  "

$ws.Range("D16").Value = "2.067.03"
$ws.Range("E16").Value = "  -1.09%  "

Set-TextCell "D17" "79.15"
$ws.Range("E17").Value = "  -3.24%  "

Set-TextCell "D18" "5.888"
$ws.Range("E18").Value = "  -5.60%  "

$ws.Range("D19").Value = "29.198.04"
$ws.Range("E19").Value = "  -0.53%  "

Set-TextCell "D20" "227.71"
$ws.Range("E20").Value = "  -0.42%  "

$ws.Range("E21").Value = "  -4.31%  "

Set-TextCell "D22" "0.9999"
$ws.Range("E22").Value = "  +0.05%  "

Set-TextCell "D23" "7.016"
$ws.Range("E23").Value = "  -4.80%  "

Set-TextCell "D24" "1.001"
$ws.Range("E24").Value = "  +0.09%  "

Set-TextCell "D25" "154.48"
$ws.Range("E25").Value = "  -2.26%  "

Set-TextCell "D26" "8.060"
$ws.Range("E26").Value = "  -5.24%  "

Set-TextCell "D27" "0.1299"
$ws.Range("E27").Value = "  -3.62%  "

Set-TextCell "D28" "16.54"
$ws.Range("E28").Value = "  -5.13%  "

Set-TextCell "D29" "1.477"
$ws.Range("E29").Value = "  +1.38%  "

Set-TextCell "D30" "0.06428"
$ws.Range("E30").Value = "  -6.51%  "

Set-TextCell "D31" "1.455"
$ws.Range("E31").Value = "  -2.14%  "

Set-TextCell "D32" "3.822"
$ws.Range("E32").Value = "  -5.42%  "

$ws.Range("E33").Value = "  -6.27%  "

$ws.Range("E34").Value = "  -1.06%  "

Set-TextCell "D35" "1.750"
$ws.Range("E35").Value = "  -3.97%  "

Set-TextCell "D36" "0.6497"
$ws.Range("E36").Value = "  -6.99%  "

Set-TextCell "D37" "2.548"
$ws.Range("E37").Value = "  -1.41%  "

Set-TextCell "D38" "2.752"
$ws.Range("E38").Value = "  -2.45%  "

$ws.Range("D39").Value = "1.214.31"
$ws.Range("E39").Value = "  -1.85%  "

Set-TextCell "D40" "0.01751"
$ws.Range("E40").Value = "  -5.13%  "

Set-TextCell "D41" "6.495"
$ws.Range("E41").Value = "  -4.53%  "

Set-TextCell "D42" "0.9321"
$ws.Range("E42").Value = "  -0.68%  "

Set-TextCell "D43" "0.9995"
$ws.Range("E43").Value = "  +0.05%  "

Set-TextCell "D44" "101.04"
$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("D45").Value = "1.983.49"
$ws.Range("E45").Value = "  -0.39%  "

Set-TextCell "D46" "63.01"
$ws.Range("E46").Value = "  -3.48%  "

$ws.Range("E47").Value = "  -0.37%  "

Set-TextCell "D48" "1.610"
$ws.Range("E48").Value = "  -5.77%  "

Set-TextCell "D49" "8.538"
$ws.Range("E49").Value = "  -4.89%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D50" "0.05528"
$ws.Range("E50").Value = "  -2.65%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell "D51" "0.1075"
$ws.Range("E51").Value = "  -5.78%  "
